# Barangay.xlsx - extend the blank "table" area on Sheet1 from row 39 down to
# row 110 (71 additional blank, but fully bordered/styled, rows spanning
# columns A:N) and move the on-screen viewport/selection down to where the
# new rows were added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing last row (39) is a fully-styled "blank" row: every cell A..N
# carries the bordered/centered/wrap-text style used throughout the table,
# with no value. Re-use that exact formatting (instead of re-building the
# border/font/alignment from scratch, which would mint new style records)
# by copying row 39's formatting across the new row block A40:N110.
$lastRow = 39
$firstNewRow = $lastRow + 1
$lastNewRow = 110

$formatSource = $ws.Range("A" + $lastRow + ":N" + $lastRow)
$newRowsRange = $ws.Range("A" + $firstNewRow + ":N" + $lastNewRow)

$formatSource.Copy()
$newRowsRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Scroll the view down to the bottom of the newly added rows and leave the
# selection / active cell parked at P107, matching where editing left off.
$ws.Range("P107").Select()
